$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks ---
# Column B (2): 24 -> ~27.11  (closest achievable through this runtime's width quantization)
$ws.Columns.Item(2).ColumnWidth = 26.33
# Column F (6) needs to become its own width bucket, split off from E:F (~36.44) to ~41.44
$ws.Columns.Item(6).ColumnWidth = 40.67

# --- Populate the new step rows (41-46) ---
# Numeric "step" markers in column A
$ws.Range("A42").Value = 0.32
$ws.Range("A43").Value = 0.33
$ws.Range("A44").Value = 0.34
$ws.Rows.Item(42).RowHeight = 28.8

# Text entries - written in the same order they were authored so the
# shared-string table comes out in the same sequence as the source edit.
$ws.Range("B41").Value = "More examples"
$ws.Range("B42").Value = "Disabling events for more optimization"
$ws.Range("C43").Value = "NPE alone logic fix"
$ws.Range("F41").Value = "NPE logic change based on brd"
$ws.Range("E42").Value = "consider updating the one cell alone, need to optimize"
$ws.Range("F42").Value = "doing some terrible hacking for the sort method, need to see if this can be done better"
$ws.Range("D42").Value = "Feature requests: print all the steps, descriptions needed for the different columns"
$ws.Range("E43").Value = "task to add the col descs"
$ws.Range("F44").Value = "investigate the activewindow FreezePane error"
$ws.Range("C44").Value = "moved steps up to have it be the first sheet"
$ws.Range("B44").Value = "optimized assets"
$ws.Range("E44").Value = "need to do the same optim for expos"
$ws.Range("B45").Value = "renamed sheet names and stuff"
$ws.Range("E45").Value = "need to write out steps"
$ws.Range("B46").Value = "freezePain optimization"
$ws.Range("E46").Value = "need to write out reasons"

# --- Selection moves to E45 ---
[void]$ws.Range("E45").Select()
